$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 44

$ws.Cells.Item($row, 1).Value = "Insert Delete GetRandom O(1)"
$ws.Cells.Item($row, 2).Value = "Class Design"
$ws.Cells.Item($row, 3).Value = "No"
$ws.Cells.Item($row, 4).Value = "Yes"
$ws.Cells.Item($row, 5).Value = "Medium"
$ws.Cells.Item($row, 6).Value = "Medium"

$ws.Cells.Item($row, 7).Value = "380 - Insert Delete GetRandom O(1)"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 7), "380 - Insert Delete GetRandom O(1)")
$ws.Cells.Item($row, 7).Style = "Hyperlink"

# Leftover formatted (but empty) cell below the new row, matching original authoring artifact
$ws.Cells.Item($row + 1, 7).Style = "Hyperlink"

# Refresh the data validation ranges to include the newly added row
$ws.Range("E2:F44").Validation.Delete()
$ws.Range("E2:F44").Validation.Add(3, 1, 1, "Easy, Medium, Hard")

$ws.Range("C2:C44").Validation.Delete()
$ws.Range("C2:C44").Validation.Add(3, 0, 1, "Yes, No")
$ws.Range("C2:C44").Validation.IgnoreBlank = $false

$ws.Range("D2:D44").Validation.Delete()
$ws.Range("D2:D44").Validation.Add(3, 1, 1, "Yes, No")

$ws.Range("B2:B44").Validation.Delete()
$ws.Range("B2:B44").Validation.Add(3, 1, 1, "Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap, Class Design")

$ws.Range("M44").Select()

$wb.Save()
